$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2358
$ws.Range("I62").Value = 2392.5
$ws.Range("K62").Value = 2392.5
$ws.Range("M62").Value = -1768.5
$ws.Range("H65").Value = 2358
$ws.Range("I65").Value = 2392.5
$ws.Range("K65").Value = 11962.5
$ws.Range("M65").Value = -8842.5
$ws.Range("H98").Value = 1165.8334
$ws.Range("I98").Value = 1165.8334
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1165.8334
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 332.1666
$ws.Range("N98").ClearContents()
$ws.Range("H109").Value = 61450
$ws.Range("J109").Value = 61450
$ws.Range("L109").Value = 61450
$ws.Range("N109").Value = -64224
$ws.Range("H122").Value = 1165.8334
$ws.Range("I122").Value = 1165.8334
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3497.5002
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1047.5002
$ws.Range("N122").ClearContents()
$ws.Range("H125").Value = 5003.1304
$ws.Range("I125").Value = 425.875
$ws.Range("J125").Value = 7444.3335
$ws.Range("K125").Value = 3832.875
$ws.Range("L125").Value = 66999.0015
$ws.Range("M125").Value = -1372.875
$ws.Range("N125").Value = -71919.0015
$ws.Range("H129").Value = 830.7547
$ws.Range("I129").Value = 281
$ws.Range("J129").Value = 914.413
$ws.Range("K129").Value = 843
$ws.Range("L129").Value = 2743.239
$ws.Range("M129").Value = 4157
$ws.Range("N129").Value = -12743.239
$ws.Range("H132").Value = 1273.209
$ws.Range("I132").Value = 1298.6936
$ws.Range("K132").Value = 3896.0808
$ws.Range("M132").Value = -1366.0808

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23343.107
$ws.Range("I32").Value = 26340.777
$ws.Range("J32").Value = 11079.909
$ws.Range("K32").Value = 26340.777
$ws.Range("L32").Value = 11079.909
$ws.Range("M32").Value = -26053.777
$ws.Range("N32").Value = -11653.909
$ws.Range("H35").Value = 20850
$ws.Range("I35").Value = 4000
$ws.Range("J35").Value = 26466.666
$ws.Range("K35").Value = 4000
$ws.Range("L35").Value = 26466.666
$ws.Range("M35").Value = -3594
$ws.Range("N35").Value = -27278.666
$ws.Range("H61").Value = 5604.847
$ws.Range("I61").Value = 3129.5518
$ws.Range("K61").Value = 3129.5518
$ws.Range("M61").Value = -2917.5518
$ws.Range("H122").Value = 31252250
$ws.Range("I122").Value = 2500
$ws.Range("K122").Value = 7500
$ws.Range("M122").Value = -5050
$ws.Range("H136").Value = 5604.847
$ws.Range("I136").Value = 3129.5518
$ws.Range("K136").Value = 9388.6554
$ws.Range("M136").Value = -6838.6554

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1314.2858
$ws.Range("I20").Value = 1350
$ws.Range("J20").Value = 1100
$ws.Range("K20").Value = 1350
$ws.Range("L20").Value = 1100
$ws.Range("M20").Value = -1103
$ws.Range("N20").Value = -1594
$ws.Range("H86").Value = 1635.1212
$ws.Range("I86").Value = 1461.8148
$ws.Range("J86").Value = 2415
$ws.Range("K86").Value = 1461.8148
$ws.Range("L86").Value = 2415
$ws.Range("M86").Value = -338.8148000000001
$ws.Range("N86").Value = -4661
$ws.Range("H89").Value = 1635.1212
$ws.Range("I89").Value = 1461.8148
$ws.Range("J89").Value = 2415
$ws.Range("K89").Value = 7309.074000000001
$ws.Range("L89").Value = 12075
$ws.Range("M89").Value = -1693.074000000001
$ws.Range("N89").Value = -23307
$ws.Range("H99").Value = 2167.889
$ws.Range("I99").Value = 2133.3333
$ws.Range("K99").Value = 2133.3333
$ws.Range("M99").Value = -635.3332999999998
$ws.Range("H134").Value = 21747.393
$ws.Range("I134").Value = 1940.0513
$ws.Range("J134").Value = 86121.25
$ws.Range("K134").Value = 5820.1539
$ws.Range("L134").Value = 258363.75
$ws.Range("M134").Value = -3285.1539
$ws.Range("N134").Value = -263433.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 161.75
$ws.Range("I7").Value = 140.625
$ws.Range("J7").Value = 175.83333
$ws.Range("K7").Value = 140.625
$ws.Range("L7").Value = 175.83333
$ws.Range("M7").Value = -27.625
$ws.Range("N7").Value = -401.83333
$ws.Range("H23").Value = 46004
$ws.Range("J23").Value = 100010
$ws.Range("L23").Value = 100010
$ws.Range("N23").Value = -100490
$ws.Range("H27").Value = 46004
$ws.Range("J27").Value = 100010
$ws.Range("L27").Value = 100010
$ws.Range("N27").Value = -100394
$ws.Range("H31").Value = 1971.375
$ws.Range("I31").Value = 1436.4286
$ws.Range("J31").Value = 5716
$ws.Range("K31").Value = 1436.4286
$ws.Range("L31").Value = 5716
$ws.Range("M31").Value = -1141.4286
$ws.Range("N31").Value = -6306
$ws.Range("H34").Value = 1971.375
$ws.Range("I34").Value = 1436.4286
$ws.Range("J34").Value = 5716
$ws.Range("K34").Value = 1436.4286
$ws.Range("L34").Value = 5716
$ws.Range("M34").Value = -1234.4286
$ws.Range("N34").Value = -6120
$ws.Range("H99").Value = 4175
$ws.Range("I99").Value = 3350
$ws.Range("K99").Value = 3350
$ws.Range("M99").Value = -1852
$ws.Range("H126").Value = 4175
$ws.Range("I126").Value = 3350
$ws.Range("K126").Value = 10050
$ws.Range("M126").Value = -7580
$ws.Range("H134").Value = 2084.6667
$ws.Range("I134").Value = 1237.2245
$ws.Range("J134").Value = 3516.5518
$ws.Range("K134").Value = 3711.6735
$ws.Range("L134").Value = 10549.6554
$ws.Range("M134").Value = -1176.6735
$ws.Range("N134").Value = -15619.6554

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1714.0358
$ws.Range("I132").Value = 2111.5557
$ws.Range("J132").Value = 1525.7368
$ws.Range("K132").Value = 19004.0013
$ws.Range("L132").Value = 13731.6312
$ws.Range("M132").Value = -16474.0013
$ws.Range("N132").Value = -18791.6312
$ws.Range("H134").Value = 4036.4062
$ws.Range("I134").Value = 3602.2856
$ws.Range("K134").Value = 10806.8568
$ws.Range("M134").Value = -5736.856800000001
$ws.Range("H136").Value = 3158.52
$ws.Range("J136").Value = 3703.15
$ws.Range("L136").Value = 11109.45
$ws.Range("N136").Value = -21309.45

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 45000
$ws.Range("J32").Value = 45000
$ws.Range("L32").Value = 45000
$ws.Range("N32").Value = -45592
$ws.Range("H122").Value = 11790.833
$ws.Range("I122").Value = 27000.5
$ws.Range("J122").Value = 4186
$ws.Range("K122").Value = 81001.5
$ws.Range("L122").Value = 12558
$ws.Range("M122").Value = -78551.5
$ws.Range("N122").Value = -17458
$ws.Range("H132").Value = 8096.485
$ws.Range("I132").Value = 5681.643
$ws.Range("J132").Value = 21619.6
$ws.Range("K132").Value = 17044.929
$ws.Range("L132").Value = 64858.8
$ws.Range("M132").Value = -14514.929
$ws.Range("N132").Value = -69918.79999999999
$ws.Range("H139").Value = 40000
$ws.Range("J139").Value = 40000
$ws.Range("L139").Value = 40000
$ws.Range("N139").Value = -50280

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 4572.6665
$ws.Range("I100").Value = 2573.25
$ws.Range("J100").Value = 6857.7144
$ws.Range("K100").Value = 2573.25
$ws.Range("L100").Value = 6857.7144
$ws.Range("M100").Value = -2032.25
$ws.Range("N100").Value = -7939.7144
$ws.Range("H122").Value = 7623.1816
$ws.Range("I122").Value = 7461.923
$ws.Range("K122").Value = 22385.769
$ws.Range("M122").Value = -19935.769

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1291.902
$ws.Range("I132").Value = 544.675
$ws.Range("K132").Value = 1634.025
$ws.Range("M132").Value = 895.9750000000001
$ws.Range("H135").Value = 142894980
$ws.Range("J135").Value = 142894980
$ws.Range("L135").Value = 142894980
$ws.Range("N135").Value = -142905120
